$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 95 (shifts rows 95-100 down to 96-101),
# mirroring the alphabetically-sorted item table gaining a new entry
# "ليفه" between "لزقه النمر بسعر القطعه" (row 94) and "محلول ملح".
$ws.Rows("95:95").Insert()

# Re-create the merged cell ranges for the newly inserted row 95
# (Insert() shifts existing merges but does not create new ones for
# the blank row it introduces).
$ws.Range("A95:B95").Merge()
$ws.Range("C95:G95").Merge()
$ws.Range("H95:K95").Merge()
$ws.Range("L95:M95").Merge()
$ws.Range("N95:O95").Merge()

# Populate the new row with the new item's data.
$ws.Cells.Item(95, 1).Value = 89
$ws.Cells.Item(95, 3).Value = "ليفه"
$ws.Cells.Item(95, 8).Value = "8:0"
$ws.Cells.Item(95, 12).Value = "0"
$ws.Cells.Item(95, 14).Value = "20.00"
$ws.Cells.Item(95, 16).Value = "20.0000"
$ws.Cells.Item(95, 17).Value = "1:0"

# Copy the styling of the neighbouring row so the new row matches the
# rest of the table's look (borders/fonts/number formats/etc.).
$ws.Range("A94:Q94").Copy()
$ws.Range("A95:Q95").PasteSpecial(-4122)

# Re-apply the values (PasteSpecial formats only, values already set above,
# but re-assert them to be safe after the paste).
$ws.Cells.Item(95, 1).Value = 89
$ws.Cells.Item(95, 3).Value = "ليفه"
$ws.Cells.Item(95, 8).Value = "8:0"
$ws.Cells.Item(95, 12).Value = "0"
$ws.Cells.Item(95, 14).Value = "20.00"
$ws.Cells.Item(95, 16).Value = "20.0000"
$ws.Cells.Item(95, 17).Value = "1:0"

# The row that used to be 98 (مناديل FINE) is now row 99; continue the
# running sequence number (it kept its old number 92 after the shift).
$ws.Cells.Item(99, 1).Value = 93

# Update the running total (old 5949.185 + new item's 20.00 = 5969.185).
# That row shifted from 99 to 100.
$ws.Cells.Item(100, 16).Value = 5969.185

# Update the generated-on timestamp in the footer (now row 101).
$ws.Cells.Item(101, 1).Value = "Wednesday, 1 October, 2025 8:45 PM"
